$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# --- Country label swaps (rows reference shared-string slots whose text
#     content got swapped between two countries), plus refreshed counts ---

# Row 36 / 37: Pakistan <-> Filipinas swapped, with refreshed counts
$ws.Range("A36").Value = "Filipinas"
Set-RowValues 36 @(3414, 168, 64, 3198, 1, 0, 152)

$ws.Range("A37").Value = "Pakistan"
Set-RowValues 37 @(3277, 120, 257, 2970, 17, 3, 50)

# Rows 206 / 208: Anguila <-> Islas Virgenes Britanicas swapped (counts unchanged)
$ws.Range("A206").Value = "Islas Virgenes Britanicas"
$ws.Range("A208").Value = "Anguila"

# Rows 209 / 210: Bonaire, San Eustaquio y Saba <-> Islas Malvinas swapped (counts unchanged)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Bonaire, San Eustaquio y Saba"

# Rows 211 / 213: Sudan del Sur <-> Timor Oriental swapped (counts unchanged)
$ws.Range("A211").Value = "Timor Oriental"
$ws.Range("A213").Value = "Sudan del Sur"

# --- Refreshed case counts for the new scrape ---

# Row 23: Australia
Set-RowValues 23 @(5788, 38, 2315, 3434, 95, 2, 39)

# Row 30: India
Set-RowValues 30 @(4298, 9, 328, 3852, 0, 0, 118)

# Row 88: Afganistan
Set-RowValues 88 @(367, 18, 17, 343, 0, 0, 7)

# Row 106: Kirguistan
Set-RowValues 106 @(216, 69, 33, 179, 5, 3, 4)

# Row 200: Belice
Set-RowValues 200 @(5, 0, 0, 4, 1, 1, 1)

# --- Timestamp footer update ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 06:52"
